# repull data, push all data, mean calculation
# Update column F (dSF) values for several rows to reflect repulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F8").Value = -3
$ws.Range("F15").Value = -3
$ws.Range("F16").Value = -5
$ws.Range("F19").Value = 0
$ws.Range("F20").Value = -2
